$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value2 = 1.02
$ws.Range("C2").Value2 = 1.028866310830283
$ws.Range("D2").Value2 = 1.03690681157448
$ws.Range("E2").Value2 = 0.992614727750844
$ws.Range("F2").Value2 = 1.045727360765864
$ws.Range("I2").Value2 = 1.033921561471159
$ws.Range("J2").Value2 = 1.034016254314096
$ws.Range("K2").Value2 = 1.039699213289268
$ws.Range("L2").Value2 = 0.9955398523335997
$ws.Range("M2").Value2 = 1.048494789185094
$ws.Range("N2").Value2 = 1.035484675998544

$ws.Range("B3").Value2 = 1.02
$ws.Range("C3").Value2 = 1.029760210865075
$ws.Range("D3").Value2 = 1.03756844474071
$ws.Range("E3").Value2 = 0.9936372048519299
$ws.Range("F3").Value2 = 1.046527995211856
$ws.Range("I3").Value2 = 1.034051854403364
$ws.Range("J3").Value2 = 1.034551229255425
$ws.Range("K3").Value2 = 1.040170937157915
$ws.Range("L3").Value2 = 0.9963617723202687
$ws.Range("M3").Value2 = 1.04910695086899
$ws.Range("N3").Value2 = 1.036020410665653

$ws.Range("B4").Value2 = 1.02
$ws.Range("C4").Value2 = 1.030339144226321
$ws.Range("D4").Value2 = 1.037996868165948
$ws.Range("E4").Value2 = 0.9942998659930998
$ws.Range("F4").Value2 = 1.047046694576457
$ws.Range("I4").Value2 = 1.034134962767797
$ws.Range("J4").Value2 = 1.034897287257762
$ws.Range("K4").Value2 = 1.040475793455666
$ws.Range("L4").Value2 = 0.9968940712668347
$ws.Range("M4").Value2 = 1.049503040600783
$ws.Range("N4").Value2 = 1.036366960110046

$ws.Range("B5").Value2 = 1.02
$ws.Range("C5").Value2 = 1.030582651035127
$ws.Range("D5").Value2 = 1.038177048131475
$ws.Range("E5").Value2 = 0.994578699834602
$ws.Range("F5").Value2 = 1.047264906090025
$ws.Range("I5").Value2 = 1.034169613645698
$ws.Range("J5").Value2 = 1.035042743365457
$ws.Range("K5").Value2 = 1.040603862698475
$ws.Range("L5").Value2 = 0.9971179600053012
$ws.Range("M5").Value2 = 1.049669550495433
$ws.Range("N5").Value2 = 1.036512622782098

$ws.Range("B6").Value2 = 1.02
$ws.Range("C6").Value2 = 1.030623544087147
$ws.Range("D6").Value2 = 1.038207305241245
$ws.Range("E6").Value2 = 0.994625531979634
$ws.Range("F6").Value2 = 1.047301553524636
$ws.Range("I6").Value2 = 1.03417541478105
$ws.Range("J6").Value2 = 1.035067164487743
$ws.Range("K6").Value2 = 1.040625360636924
$ws.Range("L6").Value2 = 0.9971555583673455
$ws.Range("M6").Value2 = 1.049697507830962
$ws.Range("N6").Value2 = 1.036537078585179

$ws.Range("B7").Value2 = 1.02
$ws.Range("C7").Value2 = 1.030342397493691
$ws.Range("D7").Value2 = 1.037999275463536
$ws.Range("E7").Value2 = 0.994303590798249
$ws.Range("F7").Value2 = 1.047049609741067
$ws.Range("I7").Value2 = 1.03413542690647
$ws.Range("J7").Value2 = 1.034899230955675
$ws.Range("K7").Value2 = 1.040477505087555
$ws.Range("L7").Value2 = 0.9968970624462089
$ws.Range("M7").Value2 = 1.049505265540487
$ws.Range("N7").Value2 = 1.036368906568233

$ws.Range("B8").Value2 = 1.02
$ws.Range("C8").Value2 = 1.029168300008079
$ws.Range("D8").Value2 = 1.037130349959317
$ws.Range("E8").Value2 = 0.9929600610674297
$ws.Range("F8").Value2 = 1.045997806291352
$ws.Range("I8").Value2 = 1.033965842643196
$ws.Range("J8").Value2 = 1.034197072962486
$ws.Range("K8").Value2 = 1.039858712642773
$ws.Range("L8").Value2 = 0.9958175282591056
$ws.Range("M8").Value2 = 1.048701675301611
$ws.Range("N8").Value2 = 1.035665751430155

$ws.Range("B9").Value2 = 1.02
$ws.Range("C9").Value2 = 1.027103433585854
$ws.Range("D9").Value2 = 1.03560158277943
$ws.Range("E9").Value2 = 0.9906006454969559
$ws.Range("F9").Value2 = 1.044149339733557
$ws.Range("I9").Value2 = 1.033657851565002
$ws.Range("J9").Value2 = 1.032959016268382
$ws.Range("K9").Value2 = 1.038765460793317
$ws.Range("L9").Value2 = 0.9939188001724441
$ws.Range("M9").Value2 = 1.04728556049517
$ws.Range("N9").Value2 = 1.034425936553538

$ws.Range("B10").Value2 = 1.02
$ws.Range("C10").Value2 = 1.025729652922502
$ws.Range("D10").Value2 = 1.034584110993785
$ws.Range("E10").Value2 = 0.989033133672735
$ws.Range("F10").Value2 = 1.042920460779829
$ws.Range("I10").Value2 = 1.033446401400987
$ws.Range("J10").Value2 = 1.032133200147101
$ws.Range("K10").Value2 = 1.038034776011988
$ws.Range("L10").Value2 = 0.9926553831429383
$ws.Range("M10").Value2 = 1.046341514226509
$ws.Range("N10").Value2 = 1.033598947678642

$ws.Range("B11").Value2 = 1.02
$ws.Range("C11").Value2 = 1.025135470412192
$ws.Range("D11").Value2 = 1.034143960033077
$ws.Range("E11").Value2 = 0.988355674866747
$ws.Range("F11").Value2 = 1.042389179166828
$ws.Range("I11").Value2 = 1.033353396694143
$ws.Range("J11").Value2 = 1.031775521553562
$ws.Range("K11").Value2 = 1.037717957902802
$ws.Range("L11").Value2 = 0.9921088820399291
$ws.Range("M11").Value2 = 1.045932757461478
$ws.Range("N11").Value2 = 1.033240761140474

$ws.Range("B12").Value2 = 1.02
$ws.Range("C12").Value2 = 1.024914867050176
$ws.Range("D12").Value2 = 1.033980533299683
$ws.Range("E12").Value2 = 0.9881042295826724
$ws.Range("F12").Value2 = 1.042191964190009
$ws.Range("I12").Value2 = 1.033318634057883
$ws.Range("J12").Value2 = 1.031642650707695
$ws.Range("K12").Value2 = 1.037600214650215
$ws.Range("L12").Value2 = 0.9919059725120875
$ws.Range("M12").Value2 = 1.045780931708344
$ws.Range("N12").Value2 = 1.033107701602766

$ws.Range("B13").Value2 = 1.02
$ws.Range("C13").Value2 = 1.024962182562235
$ws.Range("D13").Value2 = 1.034015585923625
$ws.Range("E13").Value2 = 0.9881581567098651
$ws.Range("F13").Value2 = 1.042234261710913
$ws.Range("I13").Value2 = 1.033326100550833
$ws.Range("J13").Value2 = 1.031671152529143
$ws.Range("K13").Value2 = 1.037625473804419
$ws.Range("L13").Value2 = 0.9919494934313052
$ws.Range("M13").Value2 = 1.045813498610931
$ws.Range("N13").Value2 = 1.033136243900069

$ws.Range("B14").Value2 = 1.02
$ws.Range("C14").Value2 = 1.025117233171528
$ws.Range("D14").Value2 = 1.034130449794061
$ws.Range("E14").Value2 = 0.9883348863814464
$ws.Range("F14").Value2 = 1.042372874720848
$ws.Range("I14").Value2 = 1.03335052761777
$ws.Range("J14").Value2 = 1.031764538668407
$ws.Range("K14").Value2 = 1.037708226484924
$ws.Range("L14").Value2 = 0.9920921077337197
$ws.Range("M14").Value2 = 1.045920207392742
$ws.Range("N14").Value2 = 1.033229762658362

$ws.Range("B15").Value2 = 1.02
$ws.Range("C15").Value2 = 1.025212778586953
$ws.Range("D15").Value2 = 1.034201229869684
$ws.Range("E15").Value2 = 0.9884438009545853
$ws.Range("F15").Value2 = 1.042458295610093
$ws.Range("I15").Value2 = 1.033365549263468
$ws.Range("J15").Value2 = 1.031822075250604
$ws.Range("K15").Value2 = 1.037759204844054
$ws.Range("L15").Value2 = 0.9921799884222134
$ws.Range("M15").Value2 = 1.045985954863306
$ws.Range("N15").Value2 = 1.033287380949106

$ws.Range("B16").Value2 = 1.02
$ws.Range("C16").Value2 = 1.025769101159614
$ws.Range("D16").Value2 = 1.034613331377259
$ws.Range("E16").Value2 = 0.9890781214508737
$ws.Range("F16").Value2 = 1.042955737855149
$ws.Range("I16").Value2 = 1.03345254341611
$ws.Range("J16").Value2 = 1.03215693619868
$ws.Range("K16").Value2 = 1.038055793298758
$ws.Range("L16").Value2 = 0.9926916645766087
$ws.Range("M16").Value2 = 1.046368642646055
$ws.Range("N16").Value2 = 1.033622717438137

$ws.Range("B17").Value2 = 1.02
$ws.Range("C17").Value2 = 1.026118248851345
$ws.Range("D17").Value2 = 1.03487194567232
$ws.Range("E17").Value2 = 0.989476357848556
$ws.Range("F17").Value2 = 1.043267994015217
$ws.Range("I17").Value2 = 1.033506726005514
$ws.Range("J17").Value2 = 1.032366961063054
$ws.Range("K17").Value2 = 1.038241722000354
$ws.Range("L17").Value2 = 0.9930127773699352
$ws.Range("M17").Value2 = 1.046608699382282
$ws.Range("N17").Value2 = 1.033833040561909

$ws.Range("B18").Value2 = 1.02
$ws.Range("C18").Value2 = 1.026321965697724
$ws.Range("D18").Value2 = 1.035022831526981
$ws.Range("E18").Value2 = 0.9897087662937556
$ws.Range("F18").Value2 = 1.043450207816718
$ws.Range("I18").Value2 = 1.03353819029328
$ws.Range("J18").Value2 = 1.032489455792916
$ws.Range("K18").Value2 = 1.038350129806898
$ws.Range("L18").Value2 = 0.9932001317071769
$ws.Range("M18").Value2 = 1.046748722390286
$ws.Range("N18").Value2 = 1.033955709248337

$ws.Range("B19").Value2 = 1.02
$ws.Range("C19").Value2 = 1.026391438847929
$ws.Range("D19").Value2 = 1.035074286524675
$ws.Range("E19").Value2 = 0.9897880325774034
$ws.Range("F19").Value2 = 1.043512351539784
$ws.Range("I19").Value2 = 1.03354889512225
$ws.Range("J19").Value2 = 1.032531221731466
$ws.Range("K19").Value2 = 1.038387087032041
$ws.Range("L19").Value2 = 0.9932640239640975
$ws.Range("M19").Value2 = 1.046796466907505
$ws.Range("N19").Value2 = 1.033997534499311

$ws.Range("B20").Value2 = 1.02
$ws.Range("C20").Value2 = 1.026080781885641
$ws.Range("D20").Value2 = 1.034844194603537
$ws.Range("E20").Value2 = 0.9894336180360679
$ws.Range("F20").Value2 = 1.043234483589705
$ws.Range("I20").Value2 = 1.033500927151203
$ws.Range("J20").Value2 = 1.032344428328159
$ws.Range("K20").Value2 = 1.038221777869138
$ws.Range("L20").Value2 = 0.9929783193494215
$ws.Range("M20").Value2 = 1.046582943351365
$ws.Range("N20").Value2 = 1.033810475827947

$ws.Range("B21").Value2 = 1.02
$ws.Range("C21").Value2 = 1.025071571792424
$ws.Range("D21").Value2 = 1.034096623441268
$ws.Range("E21").Value2 = 0.9882828385668249
$ws.Range("F21").Value2 = 1.042332053133044
$ws.Range("I21").Value2 = 1.033343340428271
$ws.Range("J21").Value2 = 1.031737039130138
$ws.Range("K21").Value2 = 1.037683859609657
$ws.Range("L21").Value2 = 0.9920501090198102
$ws.Range("M21").Value2 = 1.045888784174271
$ws.Range("N21").Value2 = 1.033202224067595

$ws.Range("B22").Value2 = 1.02
$ws.Range("C22").Value2 = 1.024437634895601
$ws.Range("D22").Value2 = 1.033626972279199
$ws.Range("E22").Value2 = 0.9875604150241495
$ws.Range("F22").Value2 = 1.041765393052687
$ws.Range("I22").Value2 = 1.033243006827034
$ws.Range("J22").Value2 = 1.031355074263614
$ws.Range("K22").Value2 = 1.037345285729345
$ws.Range("L22").Value2 = 0.9914670000341481
$ws.Range("M22").Value2 = 1.045452366726521
$ws.Range("N22").Value2 = 1.032819716767148

$ws.Range("B23").Value2 = 1.02
$ws.Range("C23").Value2 = 1.024773640149525
$ws.Range("D23").Value2 = 1.033875906999987
$ws.Range("E23").Value2 = 0.9879432794643023
$ws.Range("F23").Value2 = 1.042065720126906
$ws.Range("I23").Value2 = 1.033296314094231
$ws.Range("J23").Value2 = 1.031557567808822
$ws.Range("K23").Value2 = 1.037524804181474
$ws.Range("L23").Value2 = 0.991776070289318
$ws.Range("M23").Value2 = 1.045683716727217
$ws.Range("N23").Value2 = 1.033022497876418

$ws.Range("B24").Value2 = 1.02
$ws.Range("C24").Value2 = 1.026097711394272
$ws.Range("D24").Value2 = 1.034856733992016
$ws.Range("E24").Value2 = 0.9894529299347244
$ws.Range("F24").Value2 = 1.043249625261378
$ws.Range("I24").Value2 = 1.033503547834956
$ws.Range("J24").Value2 = 1.032354609929761
$ws.Range("K24").Value2 = 1.038230789891171
$ws.Range("L24").Value2 = 0.9929938892766442
$ws.Range("M24").Value2 = 1.046594581387148
$ws.Range("N24").Value2 = 1.033820671888591

$ws.Range("B25").Value2 = 1.02
$ws.Range("C25").Value2 = 1.027636763822857
$ws.Range("D25").Value2 = 1.035996512228063
$ws.Range("E25").Value2 = 0.9912096547607049
$ws.Range("F25").Value2 = 1.04462661554602
$ws.Range("I25").Value2 = 1.033738556548687
$ws.Range("J25").Value2 = 1.033279167089877
$ws.Range("K25").Value2 = 1.039048424057639
$ws.Range("L25").Value2 = 0.9944092447426414
$ws.Range("M25").Value2 = 1.047651661233394
$ws.Range("N25").Value2 = 1.034746542025921

